$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.019999999999999
$ws.Range("C2").Value = 1.008157561927207
$ws.Range("D2").Value = 1.031908695811567
$ws.Range("E2").Value = 1.010754702834615
$ws.Range("F2").Value = 1.023523560585324
$ws.Range("I2").Value = 1.031182607401444
$ws.Range("J2").Value = 1.013425020236795
$ws.Range("K2").Value = 1.034715454374724
$ws.Range("L2").Value = 1.013623937449379
$ws.Range("M2").Value = 1.02635474757224
$ws.Range("N2").Value = 1.00860477369097

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.009138964678214
$ws.Range("D3").Value = 1.03231427939092
$ws.Range("E3").Value = 1.011587648709747
$ws.Range("F3").Value = 1.024698387513373
$ws.Range("I3").Value = 1.031236483683852
$ws.Range("J3").Value = 1.014038053409991
$ws.Range("K3").Value = 1.034930779979953
$ws.Range("L3").Value = 1.014260946678546
$ws.Range("M3").Value = 1.027335468885301
$ws.Range("N3").Value = 1.008812318123289

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.009774561752514
$ws.Range("D4").Value = 1.032576308677938
$ws.Range("E4").Value = 1.012127480847482
$ws.Range("F4").Value = 1.025458468437626
$ws.Range("I4").Value = 1.031269677770116
$ws.Range("J4").Value = 1.014434748068935
$ws.Range("K4").Value = 1.035068962485437
$ws.Range("L4").Value = 1.014673354616881
$ws.Range("M4").Value = 1.027969387188669
$ws.Range("N4").Value = 1.008946484043177

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.01004190150845
$ws.Range("D5").Value = 1.032686365424027
$ws.Range("E5").Value = 1.01235463144146
$ws.Range("F5").Value = 1.025777979906885
$ws.Range("I5").Value = 1.031283232543117
$ws.Range("J5").Value = 1.014601523147745
$ws.Range("K5").Value = 1.035126778221838
$ws.Range("L5").Value = 1.014846783108316
$ws.Range("M5").Value = 1.028235724791089
$ws.Range("N5").Value = 1.009002856213794

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.010086796924366
$ws.Range("D6").Value = 1.032704838488727
$ws.Range("E6").Value = 1.01239278303239
$ws.Range("F6").Value = 1.025831625729752
$ws.Range("I6").Value = 1.031285484967648
$ws.Range("J6").Value = 1.014629525677133
$ws.Range("K6").Value = 1.035136469499866
$ws.Range("L6").Value = 1.014875905559487
$ws.Range("M6").Value = 1.028280434558232
$ws.Range("N6").Value = 1.009012319513489

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.009778133430255
$ws.Range("D7").Value = 1.032577779658626
$ws.Range("E7").Value = 1.01213051523973
$ws.Range("F7").Value = 1.025462737869963
$ws.Range("I7").Value = 1.031269860462511
$ws.Range("J7").Value = 1.014436976507994
$ws.Range("K7").Value = 1.035069736109228
$ws.Range("L7").Value = 1.014675671772493
$ws.Range("M7").Value = 1.027972946640359
$ws.Range("N7").Value = 1.008947237413718

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.008489114679871
$ws.Range("D8").Value = 1.032045848392459
$ws.Range("E8").Value = 1.011036022128824
$ws.Range("F8").Value = 1.023920621723085
$ws.Range("I8").Value = 1.031201159909769
$ws.Range("J8").Value = 1.013632193017655
$ws.Range("K8").Value = 1.034788461405408
$ws.Range("L8").Value = 1.013839171732784
$ws.Range("M8").Value = 1.026686326299016
$ws.Range("N8").Value = 1.008674940906365

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.006222028175856
$ws.Range("D9").Value = 1.03110546721852
$ws.Range("E9").Value = 1.009114012105881
$ws.Range("F9").Value = 1.021202370445667
$ws.Range("I9").Value = 1.031067366138715
$ws.Range("J9").Value = 1.012214242974637
$ws.Range("K9").Value = 1.034284094697896
$ws.Range("L9").Value = 1.012366864276665
$ws.Range("M9").Value = 1.024413987851854
$ws.Range("N9").Value = 1.008194141070506

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.004713560051089
$ws.Range("D10").Value = 1.030476620856087
$ws.Range("E10").Value = 1.007837175535151
$ws.Range("F10").Value = 1.019389633113593
$ws.Range("I10").Value = 1.030969656804606
$ws.Range("J10").Value = 1.011269088963
$ws.Range("K10").Value = 1.033942070587574
$ws.Range("L10").Value = 1.011386509442683
$ws.Range("M10").Value = 1.022895647591671
$ws.Range("N10").Value = 1.007872964236376

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.00406107179846
$ws.Range("D11").Value = 1.030203893164928
$ws.Range("E11").Value = 1.007285369078627
$ws.Range("F11").Value = 1.018604561159586
$ws.Range("I11").Value = 1.030925338441076
$ws.Range("J11").Value = 1.010859866119108
$ws.Range("K11").Value = 1.033792617693473
$ws.Range("L11").Value = 1.010962292047679
$ws.Range("M11").Value = 1.022237373567675
$ws.Range("N11").Value = 1.007733741347171

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.003818812092956
$ws.Range("D12").Value = 1.030102527034897
$ws.Range("E12").Value = 1.0070805654184
$ws.Range("F12").Value = 1.018312928463466
$ws.Range("I12").Value = 1.030908575503303
$ws.Range("J12").Value = 1.010707868123521
$ws.Range("K12").Value = 1.033736902236247
$ws.Range("L12").Value = 1.010804761759479
$ws.Range("M12").Value = 1.021992737707338
$ws.Range("N12").Value = 1.007682005264285

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.003870772938306
$ws.Range("D13").Value = 1.030124273227029
$ws.Range("E13").Value = 1.00712448916385
$ws.Range("F13").Value = 1.018375485667477
$ws.Range("I13").Value = 1.030912184824205
$ws.Range("J13").Value = 1.010740471951488
$ws.Range("K13").Value = 1.033748862508532
$ws.Range("L13").Value = 1.010838550590418
$ws.Range("M13").Value = 1.022045218521704
$ws.Range("N13").Value = 1.007693103844571

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.004041044424234
$ws.Range("D14").Value = 1.030195515480268
$ws.Range("E14").Value = 1.007268436638683
$ws.Range("F14").Value = 1.018580455160681
$ws.Range("I14").Value = 1.030923958949693
$ws.Range("J14").Value = 1.01084730180151
$ws.Range("K14").Value = 1.033788016348719
$ws.Range("L14").Value = 1.010949269666246
$ws.Range("M14").Value = 1.022217154412004
$ws.Range("N14").Value = 1.007729465284241

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.004145968020824
$ws.Range("D15").Value = 1.030239401909741
$ws.Range("E15").Value = 1.007357148894573
$ws.Range("F15").Value = 1.018706740704603
$ws.Range("I15").Value = 1.030931173498711
$ws.Range("J15").Value = 1.010913123948129
$ws.Range("K15").Value = 1.033812113598548
$ws.Range("L15").Value = 1.011017493047031
$ws.Range("M15").Value = 1.022323073413832
$ws.Range("N15").Value = 1.007751865791817

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.004756878056235
$ws.Range("D16").Value = 1.030494711953621
$ws.Range("E16").Value = 1.007873819765652
$ws.Range("F16").Value = 1.019441732681074
$ws.Range("I16").Value = 1.030972555796673
$ws.Range("J16").Value = 1.011296248501701
$ws.Range("K16").Value = 1.033951960876981
$ws.Range("L16").Value = 1.011414669350318
$ws.Range("M16").Value = 1.022939317739598
$ws.Range("N16").Value = 1.0078822008336

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.005140270126049
$ws.Range("D17").Value = 1.030654746770806
$ws.Range("E17").Value = 1.008198201335182
$ws.Range("F17").Value = 1.019902734790085
$ws.Range("I17").Value = 1.030997976358754
$ws.Range("J17").Value = 1.011536582029025
$ws.Range("K17").Value = 1.03403932185319
$ws.Range("L17").Value = 1.011663883572182
$ws.Range("M17").Value = 1.023325651286147
$ws.Range("N17").Value = 1.007963916231882

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.00536396265025
$ws.Range("D18").Value = 1.030748050444671
$ws.Range("E18").Value = 1.008387511019137
$ws.Range("F18").Value = 1.020171615575021
$ws.Range("I18").Value = 1.031012609794655
$ws.Range("J18").Value = 1.011676767872681
$ws.Range("K18").Value = 1.034090147275156
$ws.Range("L18").Value = 1.011809273258319
$ws.Range("M18").Value = 1.023550913730465
$ws.Range("N18").Value = 1.008011564838239

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.005440247283903
$ws.Range("D19").Value = 1.030779857409307
$ws.Range("E19").Value = 1.008452078245504
$ws.Range("F19").Value = 1.020263294617071
$ws.Range("I19").Value = 1.031017566497042
$ws.Range("J19").Value = 1.011724568183428
$ws.Range("K19").Value = 1.034107455207463
$ws.Range("L19").Value = 1.011858851991478
$ws.Range("M19").Value = 1.023627708938352
$ws.Range("N19").Value = 1.00802780929825

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.005099128915071
$ws.Range("D20").Value = 1.030637580870748
$ws.Range("E20").Value = 1.008163387561876
$ws.Range("F20").Value = 1.019853275070905
$ws.Range("I20").Value = 1.030995269028222
$ws.Range("J20").Value = 1.011510796165793
$ws.Range("K20").Value = 1.034029962361547
$ws.Range("L20").Value = 1.011637142423026
$ws.Range("M20").Value = 1.023284209562856
$ws.Range("N20").Value = 1.007955150451885

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.00399090085838
$ws.Range("D21").Value = 1.030174538132198
$ws.Range("E21").Value = 1.007226043221982
$ws.Range("F21").Value = 1.018520097363645
$ws.Range("I21").Value = 1.030920500069369
$ws.Range("J21").Value = 1.010815842916484
$ws.Range("K21").Value = 1.033776492083166
$ws.Range("L21").Value = 1.010916664468905
$ws.Range("M21").Value = 1.022166526988921
$ws.Range("N21").Value = 1.007718758366023

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.003294712650712
$ws.Range("D22").Value = 1.029883041784247
$ws.Range("E22").Value = 1.006637634279321
$ws.Range("F22").Value = 1.017681748104773
$ws.Range("I22").Value = 1.030871748033411
$ws.Range("J22").Value = 1.010378930722495
$ws.Range("K22").Value = 1.033615957273841
$ws.Range("L22").Value = 1.010463920055186
$ws.Range("M22").Value = 1.021463080616342
$ws.Range("N22").Value = 1.00756999918957

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.003663718415602
$ws.Range("D23").Value = 1.030037603210779
$ws.Range("E23").Value = 1.006949471962471
$ws.Range("F23").Value = 1.018126185088441
$ws.Range("I23").Value = 1.0308977572518
$ws.Range("J23").Value = 1.010610542956186
$ws.Range("K23").Value = 1.03370117004204
$ws.Range("L23").Value = 1.010703904708389
$ws.Range("M23").Value = 1.021836058614209
$ws.Range("N23").Value = 1.007648871503509

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.005117718652084
$ws.Range("D24").Value = 1.030645337531904
$ws.Range("E24").Value = 1.008179118088021
$ws.Range("F24").Value = 1.019875623831825
$ws.Range("I24").Value = 1.030996492953703
$ws.Range("J24").Value = 1.011522447677019
$ws.Range("K24").Value = 1.034034191917137
$ws.Range("L24").Value = 1.011649225513072
$ws.Range("M24").Value = 1.023302935538801
$ws.Range("N24").Value = 1.007959111375672

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.006807610107469
$ws.Range("D25").Value = 1.031348926809489
$ws.Range("E25").Value = 1.009610107520255
$ws.Range("F25").Value = 1.021905203886597
$ws.Range("I25").Value = 1.031103459131419
$ws.Range("J25").Value = 1.012580793081218
$ws.Range("K25").Value = 1.034415510570458
$ws.Range("L25").Value = 1.012747284800372
$ws.Range("M25").Value = 1.025002050543519
$ws.Range("N25").Value = 1.008318553842353
